$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a "Ship Date" column between "Sale Date" and "Items Total" in the
# two master-detail header rows (row 5 = labels, row 6 = field
# placeholders). Only those two rows gain a column, so shift the existing
# cells one column to the right by hand (copy value+format) instead of
# doing a sheet-wide column insert, which would also disturb the unrelated
# sub-range headers in rows 9-11.

# Row 5: labels ("Order No" | "Sale Date" | "Items Total" | "Amount Paid")
$ws.Range("F5").Copy($ws.Range("G5"))
$ws.Range("E5").Copy($ws.Range("F5"))
$ws.Range("E5").Value = "Ship Date"

# Row 6: field placeholders (mirrors row 5)
$ws.Range("F6").Copy($ws.Range("G6"))
$ws.Range("E6").Copy($ws.Range("F6"))
$ws.Range("E6").Value = "{{item.ShipDate}}"

# The report title in row 4 no longer spans the widened table as a single
# merged block.
$ws.Range("B4:I4").UnMerge()

# Extend the named sub-ranges so they cover the new column (F -> I).
$wb.Names.Item("Customers").RefersTo = "='Sheet 1'!`$A`$4:`$I`$13"
$wb.Names.Item("Customers_Orders").RefersTo = "='Sheet 1'!`$A`$6:`$I`$7"
$wb.Names.Item("Customers_Visitors").RefersTo = "='Sheet 1'!`$A`$10:`$I`$11"

# Match the saved selection/active cell.
$ws.Range("E6").Select() | Out-Null
